# ---------------------------------------------------------------------------
# Edit: Added min and max purity cases to computeRho.m and implemented a
# scheme for plotting the coherency matrices onto the Bloch sphere using
# Stokes parameters.
#
# Concretely, in the workbook this means:
#   1. Tiny floating point refinements to two values on the "rho_mat" sheet.
#   2. A full recomputation of the values on the "rho_min" sheet.
#   3. A brand new "rho_max" sheet (added after "rho_min") holding the
#      analogous "max purity" computation, mirroring the layout of
#      "rho_min" / "rho_mat".
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. rho_mat: two tiny numeric refinements
# ---------------------------------------------------------------------------
$wsMat = $wb.Worksheets.Item("rho_mat")
$wsMat.Range("B2").Value = 0.00049051873917709834
$wsMat.Range("B20").Value = 0.00078929553688579653

# ---------------------------------------------------------------------------
# 2. rho_min: recomputed values for every data row
# ---------------------------------------------------------------------------
$wsMin = $wb.Worksheets.Item("rho_min")

$minData = New-Object 'object[,]' 19,6
$minRows = @(
    @(0,  0.10036376636895486,   0.89963623363104506,   -0.0038412114997304313, -0.017293572703331837,  0.99999999999999978),
    @(5,  0.12538482442146476,   0.8746151755785353,    -0.13735243102516975,   -0.016409218160375137,  1),
    @(10, 0.19756687843983958,   0.80243312156016045,   -0.25702863453164704,   -0.014363176068045965,  1),
    @(15, 0.30198581433181659,   0.69801418566818352,   -0.34055184863286858,   -0.011468916211771651,  1),
    @(20, 0.43958411157855803,   0.56041588842144208,   -0.38680666172546657,   -0.0076271623921854688, 1),
    @(25, 0.57452020168969409,   0.42547979831030591,   -0.3766975967772716,    -0.0036370477396646069, 1),
    @(30, 0.69321822098959551,   0.30678177901040443,   -0.32976427663755736,   0.000097731331590671767, 1),
    @(35, 0.79448863680896575,   0.20551136319103422,   -0.24124029597275276,   0.0034784567941999762,  1),
    @(40, 0.85610120738767792,   0.14389879261232214,   -0.12285238469632448,   0.0060841180992465323,  1),
    @(45, 0.877892987347623,     0.12210701265237693,   0.0081807649160151943,  0.007570374666533034,   0.99999999999999978),
    @(50, 0.8523368181681481,    0.14766318183185181,   0.1379205118982314,     0.0077947225466580214,  0.99999999999999978),
    @(55, 0.79081315700006416,   0.20918684299993581,   0.2482170342443894,     0.0067302029385904764,  1),
    @(60, 0.68752646002804518,   0.31247353997195482,   0.33645817139793183,    0.0041738923701593721,  1),
    @(65, 0.56427406064946151,   0.43572593935053849,   0.38334176018955191,    0.00053289964899318383, 1),
    @(70, 0.43272596876307506,   0.56727403123692488,   0.38458874507582863,    -0.0038958859875596174, 1),
    @(75, 0.30447294684805898,   0.69552705315194108,   0.34197880604966485,    -0.0088966207514046096, 1),
    @(80, 0.19538615956777183,   0.80461384043222817,   0.25449263376340314,    -0.014022004246635508,  1),
    @(85, 0.12682999439279607,   0.87317000560720404,   0.14065128808143601,    -0.018467938332558959,  1),
    @(90, 0.10052981023139827,   0.89947018976860171,   -0.0019214114340316037, -0.022390069899387799,  1)
)
for ($r = 0; $r -lt 19; $r++) {
    for ($c = 0; $c -lt 6; $c++) {
        $minData[$r, $c] = $minRows[$r][$c]
    }
}
$wsMin.Range("A2:F20").Value = $minData

# Column width tweaks on rho_min
$wsMin.Columns.Item(2).ColumnWidth = 12.71
$wsMin.Columns.Item(3).ColumnWidth = 12.71
$wsMin.Columns.Item(4).ColumnWidth = 14.71
$wsMin.Columns.Item(5).ColumnWidth = 15.57
$wsMin.Columns.Item(6).ColumnWidth = 8.57

# ---------------------------------------------------------------------------
# 3. rho_max: new sheet, added right after rho_min, same layout
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsMax = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsMax.Name = "rho_max"

# Header row (reuses the existing shared strings: theta, Jxx, Jyy, beta, gamma, trace_sq)
$wsMax.Range("A1").Value = "theta"
$wsMax.Range("B1").Value = "Jxx"
$wsMax.Range("C1").Value = "Jyy"
$wsMax.Range("D1").Value = "beta"
$wsMax.Range("E1").Value = "gamma"
$wsMax.Range("F1").Value = "trace_sq"
$wsMax.Range("A1:F1").NumberFormat = "@"

$maxData = New-Object 'object[,]' 19,6
$maxRows = @(
    @(0,  0.00049051861840802521, 0.99950948138159201, -0.0048011701708678672, -0.021615416234382935, 1),
    @(5,  0.030955708681342547,   0.96904429131865755, -0.17197481007437043,   -0.020545484019007274, 1),
    @(10, 0.11925500611564953,    0.88074499388435057, -0.32358349303149325,   -0.018082369272176983, 1),
    @(15, 0.24877826587760757,    0.75122173412239246, -0.43206008540974156,   -0.014550679838937654, 1),
    @(20, 0.42285444018105367,    0.57714555981894633, -0.49391670370468899,   -0.0097391882822389225, 1),
    @(25, 0.59702770662930416,    0.40297229337069579, -0.49047242330697849,   -0.0047355534885751074, 1),
    @(30, 0.75277028594524342,    0.24722971405475658, -0.43140139720409831,   0.00012785324544165925, 1),
    @(35, 0.88677250901731708,    0.11322749098268288, -0.31683774138283183,   0.0045685004229001824,  1),
    @(40, 0.97260076801014561,    0.02739923198985433, -0.16304390480808534,   0.0080745552857337307,  0.99999999999999978),
    @(45, 0.99978264670776162,    0.00021735329223848161, 0.010819476620927537, 0.010012204550196626,  1),
    @(50, 0.96550038379422631,    0.034499616205773589, 0.18221783223087198,  0.010298232117505217,   0.99999999999999978),
    @(55, 0.88024794707427845,    0.1197520529257215,  0.32455208929998358,   0.0087999658515853889,  1),
    @(60, 0.74340730834723301,    0.25659269165276699, 0.43671905212285606,   0.005417666962836104,   1),
    @(65, 0.58267969864935854,    0.41732030135064158, 0.4931162103642871,    0.00068550177075943765, 1),
    @(70, 0.41385015739769881,    0.58614984260230119, 0.49249701921703393,   -0.0049889973657554982, 1),
    @(75, 0.25188798711515259,    0.74811201288484741, 0.43395043583561949,   -0.011289274025875064,  1),
    @(80, 0.11653015059370737,    0.8834698494062927,  0.32037366327737593,   -0.017651909214637333,  1),
    @(85, 0.032630651367381455,   0.96736934863261848, 0.17615590724660998,   -0.02312980191168083,   0.99999999999999978),
    @(90, 0.00078929562808398919, 0.99921070437191606, -0.0024011532773616806, -0.027980467258106174, 1)
)
for ($r = 0; $r -lt 19; $r++) {
    for ($c = 0; $c -lt 6; $c++) {
        $maxData[$r, $c] = $maxRows[$r][$c]
    }
}
$wsMax.Range("A2:F20").Value = $maxData

# Column widths on rho_max (mirrors the pre-edit rho_min layout)
$wsMax.Columns.Item(1).ColumnWidth = 5.86
$wsMax.Columns.Item(2).ColumnWidth = 15.71
$wsMax.Columns.Item(3).ColumnWidth = 15.71
$wsMax.Columns.Item(4).ColumnWidth = 15.43
$wsMax.Columns.Item(5).ColumnWidth = 15.71
$wsMax.Columns.Item(6).ColumnWidth = 8.57

Write-Host "rho_max sheet created and data populated."
